$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 7
$ws.Range("E1").Value = 31
$ws.Range("F1").Value = 28
$ws.Range("G1").Value = 21
$ws.Range("H1").Value = 31
$ws.Range("I1").Value = 16
$ws.Range("J1").Value = 21
$ws.Range("K1").Value = 12
$ws.Range("L1").Value = 13
$ws.Range("M1").Value = 0.001
$ws.Range("N1").Value = 0.028000000000000004
$ws.Range("O1").Value = 0.028000000000000004
$ws.Range("P1").Value = 0.027000000000000003
$ws.Range("Q1").Value = 0.031

# Column M (13) width shrinks from 7.7109375 to 5.7109375 (character units).
# The COM ColumnWidth setter here quantizes to whole screen pixels, so the
# closest reproducible value is used.
$ws.Columns.Item(13).ColumnWidth = 4.833333333333333
